# FACTURA.XLSX - update the invoice template with a new "example" client/
# issuer block and clear the invoice number / concept amount so the sheet
# goes back to a blank template (amounts recompute to 0 automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Issuer / client info block (A6:A10) ---------------------------------
$ws.Range("A6").Value  = "EJEMPLO, S.L."
$ws.Range("A7").Value  = "B12121212"
$ws.Range("A8").Value  = "C/EJEMPLO, 2 1ª PLANTA "
$ws.Range("A9").Value  = "OFICINA 1B9 - C.P. 28800"
$ws.Range("A10").Value = "EJEMPLO (MADRID)"

# --- Clear the invoice number (NÚMERO) and the CONCEPTO amount -----------
# G26 (IVA 21%) and G30 (TOTAL) are formulas driven off G22, so they fall
# back to 0 automatically once G22 is cleared.
$ws.Range("G16").ClearContents()
$ws.Range("G22").ClearContents()

# --- Update the saved view/selection state --------------------------------
$win = $excel.ActiveWindow
try { $win.ScrollRow = 4 } catch {}
try { $win.ScrollColumn = 1 } catch {}
$ws.Range("G22").Select()
